# Correct cond (skill names and img names were not corresponded)
# Swap the values of columns F (firstCond) and G (secondCond) for the
# rows whose condition/image pairing was incorrect.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(3,4,7,8,11,12,15,16,19,20,23,24,27,28,31)

foreach ($r in $rows) {
    $fCell = $ws.Range("F$r")
    $gCell = $ws.Range("G$r")
    $fVal = $fCell.Value2
    $gVal = $gCell.Value2
    $fCell.Value2 = $gVal
    $gCell.Value2 = $fVal
}

# Update the active selection to match the saved state of the workbook.
$ws.Range("G27").Select()
